# "added july 20 data"
# Updates the stat columns for Los Angeles Lakers (C), Los Angeles Clippers (E),
# Utah Jazz (I) and New Orleans Pelicans (S) across the eight stat rows (2-9),
# then clears out the now-consumed staging/scratch area (rows 12-24, cols A-I)
# that was used to enter the new game's raw numbers, and finally restores the
# last active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: PTS OFF TO ---
$ws.Range("C2").Value = 14
$ws.Range("E2").Value = 29
$ws.Range("I2").Value = 20
$ws.Range("S2").Value = 14

# --- Row 3: 2ND PTS ---
$ws.Range("C3").Value = 10
$ws.Range("E3").Value = 5
$ws.Range("I3").Value = 17
$ws.Range("S3").Value = 22

# --- Row 4: FBPS ---
$ws.Range("C4").Value = 9
$ws.Range("E4").Value = 8
$ws.Range("I4").Value = 7
$ws.Range("S4").Value = 7

# --- Row 5: PITP ---
$ws.Range("C5").Value = 40
$ws.Range("E5").Value = 24
$ws.Range("I5").Value = 56
$ws.Range("S5").Value = 48

# --- Row 6: OPP PTS OFF TO ---
$ws.Range("C6").Value = 29
$ws.Range("E6").Value = 14
$ws.Range("I6").Value = 14
$ws.Range("S6").Value = 20

# --- Row 7: OPP 2ND PTS ---
$ws.Range("C7").Value = 5
$ws.Range("E7").Value = 10
$ws.Range("I7").Value = 22
$ws.Range("S7").Value = 17

# --- Row 8: OPP FBPS ---
$ws.Range("C8").Value = 8
$ws.Range("E8").Value = 9
$ws.Range("I8").Value = 7
$ws.Range("S8").Value = 7

# --- Row 9: OPP PITP ---
$ws.Range("C9").Value = 24
$ws.Range("E9").Value = 40
$ws.Range("I9").Value = 48
$ws.Range("S9").Value = 56

# --- Clear the scratch/staging grid that was used while entering the new data ---
$ws.Range("I12:I15").Clear()
$ws.Range("F12").Copy() | Out-Null
$ws.Range("E12:E15").PasteSpecial(-4122)
$ws.Range("A16:I16").Clear()
$ws.Range("A17:D24").Clear()

# --- Restore the last selected cell ---
$ws.Range("G19").Select()
